$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = @{
    "A71" = "2025-05-11"
    "B71" = "38"
    "C71" = "37.05"
    "D71" = "0.98"
    "E71" = "0.265"
    "F71" = "0.09"
    "G71" = "5,311"
    "H71" = "7,951"
    "I71" = "8,001"
    "J71" = "7.2617"
}

foreach ($addr in $newRow.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $newRow[$addr]
    $cell.Style = "Normal"
}
